$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("list_of_char")
$ws.Activate()

# The sheet currently has a blank row at position 88 (row numbers jump from
# 87 straight to 89, i.e. row 88 has no content at all). Deleting that
# blank row shifts every row below it up by one - this both closes the gap
# and drops the now-duplicate last row (137), shrinking the used range by
# one row (A1:A137 -> A1:A136).
$ws.Rows("88").Delete()

# Leave the view the way it ends up after that delete: cell A88 (now
# holding what used to be row 89's value) selected/active, scrolled down
# so row 115 is near the top of the window.
$ws.Range("A88").Select()
$excel.ActiveWindow.ScrollRow = 115
